$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "FFR_LF"
$ws.Range("C1").Value = "LF_CA"

$ws.Range("B2").Value = 1.271001727713975
$ws.Range("C2").Value = 0.4195935936794408

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.00001685781251703489
